$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook-level: rename the dated "FOHM" sheet and move the active tab.
# ---------------------------------------------------------------------------
$shFohm = $wb.Sheets.Item("FOHM  5 Jun 2020")
$shFohm.Name = "FOHM  6 Jun 2020"

$shDeaths = $wb.Sheets.Item("Antal avlidna per dag")
$shDeaths.Activate()

# ---------------------------------------------------------------------------
# Sheet 1: "Antal per dag region"
# ---------------------------------------------------------------------------
$ws1 = $wb.Sheets.Item("Antal per dag region")

# Row 105 (2020-06-01): N 59->58, P 6->7
$ws1.Cells.Item(105, 14).Value = 58
$ws1.Cells.Item(105, 16).Value = 7

# Row 123 (2020-06-05): B 997->1042, P 57->58, S 25->27, U 363->405
$ws1.Cells.Item(123, 2).Value = 1042
$ws1.Cells.Item(123, 16).Value = 58
$ws1.Cells.Item(123, 19).Value = 27
$ws1.Cells.Item(123, 21).Value = 405

# Row 124 (2020-06-06): full row of data replaced
$row124 = @(43987, 1016, 10, 20, 1, 43, 21, 13, 42, 14, 19, 14, 34, 235, 3, 34, 8, 5, 21, 19, 373, 46, 41)
for ($i = 0; $i -lt $row124.Length; $i++) {
    $ws1.Cells.Item(124, $i + 1).Value = $row124[$i]
}

# Row 125 (2020-06-07): brand-new row - copy formatting down from row 124 first
$ws1.Range("A124:W124").Copy()
$ws1.Range("A125:W125").PasteSpecial(-4122)
$row125 = @(43988, 114, 0, 0, 0, 1, 1, 11, 9, 1, 0, 1, 1, 27, 0, 1, 3, 0, 1, 0, 27, 0, 30)
for ($i = 0; $i -lt $row125.Length; $i++) {
    $ws1.Cells.Item(125, $i + 1).Value = $row125[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Antal avlidna per dag"
# ---------------------------------------------------------------------------
$ws2 = $wb.Sheets.Item("Antal avlidna per dag")

$ws2.Cells.Item(80, 2).Value = 37
$ws2.Cells.Item(81, 2).Value = 33
$ws2.Cells.Item(82, 2).Value = 28
$ws2.Cells.Item(83, 2).Value = 38
$ws2.Cells.Item(84, 2).Value = 28
$ws2.Cells.Item(85, 2).Value = 19
$ws2.Cells.Item(87, 2).Value = 10

# Row 88 used to hold the "Uppgift saknas" aggregate; it now becomes a dated
# row (2020-06-06) like the others, copying the date-cell formatting from
# row 87 first so the number format/border match.
$ws2.Range("A87:B87").Copy()
$ws2.Range("A88:B88").PasteSpecial(-4122)
$ws2.Range("A88").Value = 43987
$ws2.Range("B88").Value = 2

# New row 89: the "Uppgift saknaa" aggregate row, copying formatting from row 88.
$ws2.Range("A88:B88").Copy()
$ws2.Range("A89:B89").PasteSpecial(-4122)
$ws2.Range("A89").Value = "Uppgift saknaa"
$ws2.Range("B89").Value = 10
$ws2.Range("B89").Select()

# ---------------------------------------------------------------------------
# Sheet 3: "Antal intensivvårdade per dag"
# ---------------------------------------------------------------------------
$ws3 = $wb.Sheets.Item("Antal intensivvårdade per dag")

$ws3.Cells.Item(35, 2).Value = 47
$ws3.Cells.Item(55, 2).Value = 33
$ws3.Cells.Item(76, 2).Value = 13
$ws3.Cells.Item(80, 2).Value = 16
$ws3.Cells.Item(85, 2).Value = 19
$ws3.Cells.Item(92, 2).Value = 11
$ws3.Cells.Item(93, 2).Value = 4

# ---------------------------------------------------------------------------
# Sheet 4: "Totalt antal per region"
# ---------------------------------------------------------------------------
$ws4 = $wb.Sheets.Item("Totalt antal per region")

$ws4.Cells.Item(2, 2).Value = 274
$ws4.Cells.Item(2, 3).Value = 171.67274475097656

$ws4.Cells.Item(3, 2).Value = 1297
$ws4.Cells.Item(3, 3).Value = 450.400390625
$ws4.Cells.Item(3, 4).Value = 61

$ws4.Cells.Item(4, 2).Value = 89
$ws4.Cells.Item(4, 3).Value = 149.11369323730469

$ws4.Cells.Item(5, 2).Value = 1470
$ws4.Cells.Item(5, 3).Value = 511.5142822265625
$ws4.Cells.Item(5, 5).Value = 112

$ws4.Cells.Item(6, 2).Value = 887
$ws4.Cells.Item(6, 3).Value = 265.68978881835938
$ws4.Cells.Item(6, 4).Value = 31
$ws4.Cells.Item(6, 5).Value = 57

$ws4.Cells.Item(7, 2).Value = 792
$ws4.Cells.Item(7, 3).Value = 605.45831298828125

$ws4.Cells.Item(8, 2).Value = 1719
$ws4.Cells.Item(8, 3).Value = 472.77359008789063

$ws4.Cells.Item(9, 2).Value = 405
$ws4.Cells.Item(9, 3).Value = 165.0057373046875

$ws4.Cells.Item(10, 2).Value = 916
$ws4.Cells.Item(10, 3).Value = 454.6605224609375

$ws4.Cells.Item(11, 2).Value = 495
$ws4.Cells.Item(11, 3).Value = 197.92637634277344

$ws4.Cells.Item(12, 2).Value = 1957
$ws4.Cells.Item(12, 3).Value = 142.03524780273438
$ws4.Cells.Item(12, 4).Value = 97
$ws4.Cells.Item(12, 5).Value = 200

$ws4.Cells.Item(13, 2).Value = 14571
$ws4.Cells.Item(13, 3).Value = 612.97869873046875
$ws4.Cells.Item(13, 4).Value = 822
$ws4.Cells.Item(13, 5).Value = 2137

$ws4.Cells.Item(14, 2).Value = 1657
$ws4.Cells.Item(14, 3).Value = 556.89990234375

$ws4.Cells.Item(15, 2).Value = 2171
$ws4.Cells.Item(15, 3).Value = 565.7874755859375

$ws4.Cells.Item(16, 2).Value = 606
$ws4.Cells.Item(16, 3).Value = 214.57859802246094

$ws4.Cells.Item(17, 2).Value = 523
$ws4.Cells.Item(17, 3).Value = 192.46621704101563

$ws4.Cells.Item(18, 2).Value = 853
$ws4.Cells.Item(18, 3).Value = 347.67083740234375

$ws4.Cells.Item(19, 2).Value = 1436
$ws4.Cells.Item(19, 3).Value = 520.58221435546875

$ws4.Cells.Item(20, 2).Value = 7589
$ws4.Cells.Item(20, 3).Value = 439.71746826171875
$ws4.Cells.Item(20, 4).Value = 349
$ws4.Cells.Item(20, 5).Value = 602

$ws4.Cells.Item(21, 2).Value = 1942
$ws4.Cells.Item(21, 3).Value = 637.128662109375

$ws4.Cells.Item(22, 2).Value = 2238
$ws4.Cells.Item(22, 3).Value = 480.77853393554688
$ws4.Cells.Item(22, 5).Value = 198

# ---------------------------------------------------------------------------
# Sheet 5: "Totalt antal per kön"
# ---------------------------------------------------------------------------
$ws5 = $wb.Sheets.Item("Totalt antal per kön")

$ws5.Cells.Item(2, 2).Value = 17493
$ws5.Cells.Item(2, 3).Value = 1595
$ws5.Cells.Item(2, 4).Value = 2554

$ws5.Cells.Item(3, 2).Value = 26393
$ws5.Cells.Item(3, 3).Value = 567
$ws5.Cells.Item(3, 4).Value = 2102

# ---------------------------------------------------------------------------
# Sheet 6: "Totalt antal per åldersgrupp"
# ---------------------------------------------------------------------------
$ws6 = $wb.Sheets.Item("Totalt antal per åldersgrupp")

$ws6.Cells.Item(2, 2).Value = 209
$ws6.Cells.Item(2, 3).Value = 5

$ws6.Cells.Item(3, 2).Value = 580
$ws6.Cells.Item(3, 3).Value = 9

$ws6.Cells.Item(4, 2).Value = 4343
$ws6.Cells.Item(4, 3).Value = 81

$ws6.Cells.Item(5, 2).Value = 5806
$ws6.Cells.Item(5, 3).Value = 96

$ws6.Cells.Item(6, 2).Value = 6852
$ws6.Cells.Item(6, 3).Value = 246

$ws6.Cells.Item(7, 2).Value = 8121
$ws6.Cells.Item(7, 3).Value = 568
$ws6.Cells.Item(7, 4).Value = 137

$ws6.Cells.Item(8, 2).Value = 5218
$ws6.Cells.Item(8, 3).Value = 649

$ws6.Cells.Item(9, 2).Value = 4391
$ws6.Cells.Item(9, 3).Value = 421
$ws6.Cells.Item(9, 4).Value = 1022

$ws6.Cells.Item(10, 2).Value = 5396
$ws6.Cells.Item(10, 3).Value = 86
$ws6.Cells.Item(10, 4).Value = 1914

$ws6.Cells.Item(11, 2).Value = 2957
$ws6.Cells.Item(11, 4).Value = 1184

Write-Output "edit complete"
